$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells
$ws.Range("B2").Value = 106
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 101

# Add new rows 4-6, copying the style of A3 (style index 1) for column A
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 38

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 35

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 35

# Apply the same style as A3 to the new A4:A6 cells
$ws.Range("A3").Copy()
$ws.Range("A4:A6").PasteSpecial(-4122)
